$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New DRS review rows (Match 53 & 54) being appended to the data table,
# mirroring columns A:Q -> Match, Home Team, Away Team, Innings, Batting Team,
# Fielding Team, Over, Review By, Umpire, Umpire Abbreviation,
# Decision Challenged, Original Decision, DRS Decision, Batter, Bowler,
# Result, Umpires Call
$newRows = @(
    @(53,"PBKS","CSK",1,"CSK","PBKS",9,"CSK","A Nand Kishore","ANK","Wicket","Out","Out","DJ Mitchell","HV Patel","Unsuccessful","Yes"),
    @(53,"PBKS","CSK",2,"PBKS","CSK",19,"CSK","A Nand Kishore","ANK","Wide","Called","Not Called","K Rabada","RJ Gleeson","Successful","No"),
    @(54,"LSG","KKR",1,"KKR","LSG",15,"LSG","YC Barde","YCB","Wide","Called","Called","AD Russell","Naveen-ul-Haq","Unsuccessful","No"),
    @(54,"LSG","KKR",1,"KKR","LSG",17,"KKR","YC Barde","YCB","Wide","Not Called","Called","SS Iyer","Yash Thakur","Successful","No"),
    @(54,"LSG","KKR",1,"KKR","LSG",20,"KKR","MV Saidharshan Kumar","MVSK","Wide","Not Called","Called","SS Iyer","Yash Thakur","Successful","No"),
    @(54,"LSG","KKR",2,"LSG","KKR",3,"LSG","MV Saidharshan Kumar","MVSK","Wicket","Out","Not Out","KL Rahul","VG Arora","Successful","No"),
    @(54,"LSG","KKR",2,"LSG","KKR",9,"LSG","MV Saidharshan Kumar","MVSK","Wicket","Out","Out","DJ Hooda","CV Varun","Unsuccessful","Yes"),
    @(54,"LSG","KKR",2,"LSG","KKR",17,"LSG","MV Saidharshan Kumar","MVSK","Wicket","Out","Out","Ravi Bishnoi","Harshit Rana","Unsuccessful","Yes")
)

$startRow = 220
$r = $startRow
foreach ($row in $newRows) {
    for ($i = 0; $i -lt $row.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $row[$i]
    }
    $r++
}

# Match the sheet view state left behind by the author's last save:
# scrolled so row 202 / column C is the top-left visible cell, with N223
# as the active selection.
$excel.ActiveWindow.ScrollRow = 202
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("N223").Select()
